# This script updates the NATMI LR-pairs data table (Pgf-Nrp2) on the active
# worksheet with recomputed values reflecting the new TPM-based calculations.
# All edits are literal value replacements in columns E-T for rows 2-17,
# matching the updated output of the NATMI analysis script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 5.961753666666667
$ws.Cells.Item(2, 8).Value2 = 17.885261
$ws.Cells.Item(2, 9).Value2 = 0.7691652453336842
$ws.Cells.Item(2, 10).Value2 = 0.7691652453336842
$ws.Cells.Item(2, 13).Value2 = 57.65261933333333
$ws.Cells.Item(2, 14).Value2 = 172.957858
$ws.Cells.Item(2, 15).Value2 = 0.6817060950001529
$ws.Cells.Item(2, 16).Value2 = 0.6817060950001529
$ws.Cells.Item(2, 17).Value2 = 343.7107147034375
$ws.Cells.Item(2, 18).Value2 = 3093.396432330938
$ws.Cells.Item(2, 19).Value2 = 0.5243446358062605
$ws.Cells.Item(2, 20).Value2 = 0.5243446358062605
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 5.961753666666667
$ws.Cells.Item(3, 8).Value2 = 17.885261
$ws.Cells.Item(3, 9).Value2 = 0.7691652453336842
$ws.Cells.Item(3, 10).Value2 = 0.7691652453336842
$ws.Cells.Item(3, 15).Value2 = 0.1019529789289588
$ws.Cells.Item(3, 16).Value2 = 0.1019529789289588
$ws.Cells.Item(3, 17).Value2 = 51.40386966000221
$ws.Cells.Item(3, 18).Value2 = 462.6348269400199
$ws.Cells.Item(3, 19).Value2 = 0.07841868805039255
$ws.Cells.Item(3, 20).Value2 = 0.07841868805039255
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 5.961753666666667
$ws.Cells.Item(4, 8).Value2 = 17.885261
$ws.Cells.Item(4, 9).Value2 = 0.7691652453336842
$ws.Cells.Item(4, 10).Value2 = 0.7691652453336842
$ws.Cells.Item(4, 13).Value2 = 2.790736
$ws.Cells.Item(4, 14).Value2 = 8.372208000000001
$ws.Cells.Item(4, 15).Value2 = 0.0329987043561157
$ws.Cells.Item(4, 16).Value2 = 0.0329987043561157
$ws.Cells.Item(4, 17).Value2 = 16.63768058069867
$ws.Cells.Item(4, 18).Value2 = 149.739125226288
$ws.Cells.Item(4, 19).Value2 = 0.02538145653176545
$ws.Cells.Item(4, 20).Value2 = 0.02538145653176544
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 5.961753666666667
$ws.Cells.Item(5, 8).Value2 = 17.885261
$ws.Cells.Item(5, 9).Value2 = 0.7691652453336842
$ws.Cells.Item(5, 10).Value2 = 0.7691652453336842
$ws.Cells.Item(5, 13).Value2 = 15.50544933333333
$ws.Cells.Item(5, 14).Value2 = 46.516348
$ws.Cells.Item(5, 15).Value2 = 0.1833422217147727
$ws.Cells.Item(5, 16).Value2 = 0.1833422217147727
$ws.Cells.Item(5, 17).Value2 = 92.43966941631422
$ws.Cells.Item(5, 18).Value2 = 831.957024746828
$ws.Cells.Item(5, 19).Value2 = 0.1410204649452658
$ws.Cells.Item(5, 20).Value2 = 0.1410204649452658
$ws.Cells.Item(6, 9).Value2 = 0.07350643921898506
$ws.Cells.Item(6, 10).Value2 = 0.07350643921898506
$ws.Cells.Item(6, 13).Value2 = 57.65261933333333
$ws.Cells.Item(6, 14).Value2 = 172.957858
$ws.Cells.Item(6, 15).Value2 = 0.6817060950001529
$ws.Cells.Item(6, 16).Value2 = 0.6817060950001529
$ws.Cells.Item(6, 17).Value2 = 32.84723394945066
$ws.Cells.Item(6, 18).Value2 = 295.625105545056
$ws.Cells.Item(6, 19).Value2 = 0.05010978763734039
$ws.Cells.Item(6, 20).Value2 = 0.05010978763734039
$ws.Cells.Item(7, 9).Value2 = 0.07350643921898506
$ws.Cells.Item(7, 10).Value2 = 0.07350643921898506
$ws.Cells.Item(7, 15).Value2 = 0.1019529789289588
$ws.Cells.Item(7, 16).Value2 = 0.1019529789289588
$ws.Cells.Item(7, 19).Value2 = 0.007494200448835977
$ws.Cells.Item(7, 20).Value2 = 0.007494200448835977
$ws.Cells.Item(8, 9).Value2 = 0.07350643921898506
$ws.Cells.Item(8, 10).Value2 = 0.07350643921898506
$ws.Cells.Item(8, 13).Value2 = 2.790736
$ws.Cells.Item(8, 14).Value2 = 8.372208000000001
$ws.Cells.Item(8, 15).Value2 = 0.0329987043561157
$ws.Cells.Item(8, 16).Value2 = 0.0329987043561157
$ws.Cells.Item(8, 17).Value2 = 1.590005091584
$ws.Cells.Item(8, 18).Value2 = 14.310045824256
$ws.Cells.Item(8, 19).Value2 = 0.002425617256058076
$ws.Cells.Item(8, 20).Value2 = 0.002425617256058076
$ws.Cells.Item(9, 9).Value2 = 0.07350643921898506
$ws.Cells.Item(9, 10).Value2 = 0.07350643921898506
$ws.Cells.Item(9, 13).Value2 = 15.50544933333333
$ws.Cells.Item(9, 14).Value2 = 46.516348
$ws.Cells.Item(9, 15).Value2 = 0.1833422217147727
$ws.Cells.Item(9, 16).Value2 = 0.1833422217147727
$ws.Cells.Item(9, 17).Value2 = 8.834136724970667
$ws.Cells.Item(9, 18).Value2 = 79.507230524736
$ws.Cells.Item(9, 19).Value2 = 0.01347683387675062
$ws.Cells.Item(9, 20).Value2 = 0.01347683387675062
$ws.Cells.Item(10, 7).Value2 = 1.143196
$ws.Cells.Item(10, 8).Value2 = 3.429588
$ws.Cells.Item(10, 9).Value2 = 0.1474912720263607
$ws.Cells.Item(10, 10).Value2 = 0.1474912720263607
$ws.Cells.Item(10, 13).Value2 = 57.65261933333333
$ws.Cells.Item(10, 14).Value2 = 172.957858
$ws.Cells.Item(10, 15).Value2 = 0.6817060950001529
$ws.Cells.Item(10, 16).Value2 = 0.6817060950001529
$ws.Cells.Item(10, 17).Value2 = 65.90824381138931
$ws.Cells.Item(10, 18).Value2 = 593.1741943025039
$ws.Cells.Item(10, 19).Value2 = 0.1005456990996956
$ws.Cells.Item(10, 20).Value2 = 0.1005456990996956
$ws.Cells.Item(11, 7).Value2 = 1.143196
$ws.Cells.Item(11, 8).Value2 = 3.429588
$ws.Cells.Item(11, 9).Value2 = 0.1474912720263607
$ws.Cells.Item(11, 10).Value2 = 0.1474912720263607
$ws.Cells.Item(11, 15).Value2 = 0.1019529789289588
$ws.Cells.Item(11, 16).Value2 = 0.1019529789289588
$ws.Cells.Item(11, 17).Value2 = 9.856948385573331
$ws.Cells.Item(11, 18).Value2 = 88.71253547015999
$ws.Cells.Item(11, 19).Value2 = 0.01503717454910888
$ws.Cells.Item(11, 20).Value2 = 0.01503717454910888
$ws.Cells.Item(12, 7).Value2 = 1.143196
$ws.Cells.Item(12, 8).Value2 = 3.429588
$ws.Cells.Item(12, 9).Value2 = 0.1474912720263607
$ws.Cells.Item(12, 10).Value2 = 0.1474912720263607
$ws.Cells.Item(12, 13).Value2 = 2.790736
$ws.Cells.Item(12, 14).Value2 = 8.372208000000001
$ws.Cells.Item(12, 15).Value2 = 0.0329987043561157
$ws.Cells.Item(12, 16).Value2 = 0.0329987043561157
$ws.Cells.Item(12, 17).Value2 = 3.190358232256
$ws.Cells.Item(12, 18).Value2 = 28.713224090304
$ws.Cells.Item(12, 19).Value2 = 0.004867020880705314
$ws.Cells.Item(12, 20).Value2 = 0.004867020880705313
$ws.Cells.Item(13, 7).Value2 = 1.143196
$ws.Cells.Item(13, 8).Value2 = 3.429588
$ws.Cells.Item(13, 9).Value2 = 0.1474912720263607
$ws.Cells.Item(13, 10).Value2 = 0.1474912720263607
$ws.Cells.Item(13, 13).Value2 = 15.50544933333333
$ws.Cells.Item(13, 14).Value2 = 46.516348
$ws.Cells.Item(13, 15).Value2 = 0.1833422217147727
$ws.Cells.Item(13, 16).Value2 = 0.1833422217147727
$ws.Cells.Item(13, 17).Value2 = 17.72576765606933
$ws.Cells.Item(13, 18).Value2 = 159.531908904624
$ws.Cells.Item(13, 19).Value2 = 0.02704137749685087
$ws.Cells.Item(13, 20).Value2 = 0.02704137749685087
$ws.Cells.Item(14, 7).Value2 = 0.07624633333333333
$ws.Cells.Item(14, 8).Value2 = 0.228739
$ws.Cells.Item(14, 9).Value2 = 0.009837043420970016
$ws.Cells.Item(14, 10).Value2 = 0.009837043420970016
$ws.Cells.Item(14, 13).Value2 = 57.65261933333333
$ws.Cells.Item(14, 14).Value2 = 172.957858
$ws.Cells.Item(14, 15).Value2 = 0.6817060950001529
$ws.Cells.Item(14, 16).Value2 = 0.6817060950001529
$ws.Cells.Item(14, 17).Value2 = 4.395800831229111
$ws.Cells.Item(14, 18).Value2 = 39.562207481062
$ws.Cells.Item(14, 19).Value2 = 0.006705972456856415
$ws.Cells.Item(14, 20).Value2 = 0.006705972456856415
$ws.Cells.Item(15, 7).Value2 = 0.07624633333333333
$ws.Cells.Item(15, 8).Value2 = 0.228739
$ws.Cells.Item(15, 9).Value2 = 0.009837043420970016
$ws.Cells.Item(15, 10).Value2 = 0.009837043420970016
$ws.Cells.Item(15, 15).Value2 = 0.1019529789289588
$ws.Cells.Item(15, 16).Value2 = 0.1019529789289588
$ws.Cells.Item(15, 17).Value2 = 0.6574167266644444
$ws.Cells.Item(15, 18).Value2 = 5.916750539979999
$ws.Cells.Item(15, 19).Value2 = 0.001002915880621409
$ws.Cells.Item(15, 20).Value2 = 0.001002915880621409
$ws.Cells.Item(16, 7).Value2 = 0.07624633333333333
$ws.Cells.Item(16, 8).Value2 = 0.228739
$ws.Cells.Item(16, 9).Value2 = 0.009837043420970016
$ws.Cells.Item(16, 10).Value2 = 0.009837043420970016
$ws.Cells.Item(16, 13).Value2 = 2.790736
$ws.Cells.Item(16, 14).Value2 = 8.372208000000001
$ws.Cells.Item(16, 15).Value2 = 0.0329987043561157
$ws.Cells.Item(16, 16).Value2 = 0.0329987043561157
$ws.Cells.Item(16, 17).Value2 = 0.2127833873013334
$ws.Cells.Item(16, 18).Value2 = 1.915050485712
$ws.Cells.Item(16, 19).Value2 = 0.0003246096875868626
$ws.Cells.Item(16, 20).Value2 = 0.0003246096875868625
$ws.Cells.Item(17, 7).Value2 = 0.07624633333333333
$ws.Cells.Item(17, 8).Value2 = 0.228739
$ws.Cells.Item(17, 9).Value2 = 0.009837043420970016
$ws.Cells.Item(17, 10).Value2 = 0.009837043420970016
$ws.Cells.Item(17, 13).Value2 = 15.50544933333333
$ws.Cells.Item(17, 14).Value2 = 46.516348
$ws.Cells.Item(17, 15).Value2 = 0.1833422217147727
$ws.Cells.Item(17, 16).Value2 = 0.1833422217147727
$ws.Cells.Item(17, 17).Value2 = 1.182233658352444
$ws.Cells.Item(17, 18).Value2 = 10.640102925172
$ws.Cells.Item(17, 19).Value2 = 0.007494200448835977
$ws.Cells.Item(17, 20).Value2 = 0.007494200448835977

Write-Host "Updated cell values with new TPM-based NATMI results"
